# The commit duplicates the last paragraph's original text
# ("哦哦拉拉拉♪(^∇^*)") into a new paragraph (followed by a blank
# paragraph), and then changes the wording of the original last
# paragraph (which keeps the _GoBack bookmark) to "我是美丽小天使".
#
# We rebuild the target paragraph plus its two new predecessors in one
# shot via Range.InsertXML so the resulting w:pPr/w:rPr "hint" values
# (eastAsia vs default) and the empty trailing paragraph match exactly
# what Word produced when a user typed two Enters in front of that
# text and then edited the wording.

$d = $word.ActiveDocument

$oldText = "哦哦拉拉拉♪(^∇^*)"
$newText = "我是美丽小天使"

# Locate the paragraph that currently holds the old text (it's the
# last paragraph in the body, carrying the _GoBack bookmark) without
# hard-coding its index.
$target = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs($i)
    $text = $para.Range.Text
    if ($text.TrimEnd([char]13) -eq $oldText) {
        $target = $para
    }
}

if ($target -eq $null) {
    throw "could not find paragraph containing '$oldText'"
}

$w = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

$xml = "<w:p xmlns:w='$w'>" +
         "<w:pPr><w:rPr><w:rFonts w:hint='eastAsia'/><w:lang w:val='en-US' w:eastAsia='zh-CN'/></w:rPr></w:pPr>" +
         "<w:r><w:rPr><w:rFonts w:hint='eastAsia'/><w:lang w:val='en-US' w:eastAsia='zh-CN'/></w:rPr><w:t>$oldText</w:t></w:r>" +
       "</w:p>" +
       "<w:p xmlns:w='$w'>" +
         "<w:pPr><w:rPr><w:rFonts w:hint='eastAsia'/><w:lang w:val='en-US' w:eastAsia='zh-CN'/></w:rPr></w:pPr>" +
       "</w:p>" +
       "<w:p xmlns:w='$w'>" +
         "<w:pPr><w:rPr><w:rFonts w:hint='default'/><w:lang w:val='en-US' w:eastAsia='zh-CN'/></w:rPr></w:pPr>" +
         "<w:r><w:rPr><w:rFonts w:hint='eastAsia'/><w:lang w:val='en-US' w:eastAsia='zh-CN'/></w:rPr><w:t>$newText</w:t></w:r>" +
         "<w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/>" +
       "</w:p>"

$target.Range.InsertXML($xml)

Write-Output "done"
